# Apply the "Drop from NY" column addition to the UON24 Entries sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H2 - matches the formatting already used by row 2's
# custom row format (bold 9pt Arial, wrapped, rotated 180).
$ws.Range("H2").Value = "Drop from NY (first five based on pedigree after that ??)"
$h2 = $ws.Range("H2")
$h2.Font.Name = "Arial"
$h2.Font.Size = 9
$h2.Font.Bold = $true
$h2.WrapText = $true
$h2.Orientation = 180

# New H column data values (bold 9pt Arial, wrapped - matches the body
# rows' custom row format used throughout columns A:G).
$hCells = "H16", "H19", "H26", "H27", "H29", "H31", "H33", "H35"
$hValues = @{ "H16" = 8; "H19" = 1; "H26" = 2; "H27" = 3; "H29" = 4; "H31" = 6; "H33" = 5; "H35" = 7 }

foreach ($addr in $hCells) {
    $cell = $ws.Range($addr)
    $cell.Value = $hValues[$addr]
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 9
    $cell.Font.Bold = $true
    $cell.WrapText = $true
}

# Update selection to match the recorded end-state of the edit session
$ws.Range("H22").Select() | Out-Null
